$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Build unicode punctuation we need (en dash, curly quotes) ---
$endash = [char]0x2013
$ldq    = [char]0x201C
$rdq    = [char]0x201D

# --- Row 34: trim reviewer-response text, update "changes to thesis" text ---
$ws.Range("B34").Value = "Chemistry is not my strongest suit, so I do not pick up errors as readily as I should " + $endash + " thanks for pointing these out. "
$ws.Range("C34").Value = "Low wavelengths discussed in equation set 1.2 are updated to 350nm, updated sentence dealing with M abundance."

# --- Row 37: add response + changes-to-thesis text for the VOC-sources comment ---
$ws.Range("B37").Value = "Thanks, it is good to have feedback on what may be missing from the intro. Instead of shuffling pyro and anthro emissions into the intro, I hope it is ok to just mention that they are important and point to the sections where I deal with them specifically."
$ws.Range("C37").Value = "Added to third paragraph in 1.3: " + $ldq + "Other major emission sources of VOC (anthropogenic and pyrogenic) are also important." + $rdq + " and " + $ldq + "This thesis mostly focuses on biogenic emissions, with influences from pyrogenic and anthropogenic emissions removed (Section 2.7)" + $rdq + ". Added to first paragraph in 1.3: " + $ldq + "These properties are largely dictated by the chemical makeup of the individual compounds. A compound's atmospheric lifetime is strongly related to its reactivity (and the concentration of reactants), with more reactive compounds having shorter atmospheric lifetimes." + $rdq

# --- Row 40: add response + changes-to-thesis text for the biogenic-isoprene comment ---
$ws.Range("B40").Value = "Fair point."
$ws.Range("C40").Value = "Added Section 1.3.2: Biogenic emissions modelling, drawing sentences from several sections in chapter 1 and 3"

# --- Row heights grow because the new wrapped text needs more vertical space ---
$ws.Range("A37").RowHeight = 157.45
$ws.Range("A40").RowHeight = 171.6

# --- Selection / view state matches where the author ended up editing ---
$ws.Range("C40").Select()
